# Conduct S3 tests for Graphql with resolved fields
#
# Adds a per-row "average of the 5 runs" column (L) to each of the three
# result sheets, normalises a redundant duplicate cell style that Excel
# folds away when the workbook is touched, and restores each sheet's
# on-screen selection / scroll position.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "1 rekord"
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)

# The last-column (K) header/value cells in rows 8:18 carried a style
# that duplicated style index 1 (border + centred alignment) except for
# a no-op "applyFill". Re-applying the same alignment collapses the
# cell onto the existing, equivalent style.
$ws1.Range("K8:K18").HorizontalAlignment = -4108
$ws1.Range("K8:K18").VerticalAlignment = -4108

# New column L: average of the five run columns (G:K) for every data row.
$ws1.Range("L14").Formula = "=AVERAGE(G14:K14)"
$ws1.Range("L15:L38").Formula = "=AVERAGE(G15:K15)"

# ---------------------------------------------------------------------
# Sheet "100 rekordów"
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)

$ws2.Range("J8:K18").HorizontalAlignment = -4108
$ws2.Range("J8:K18").VerticalAlignment = -4108

$ws2.Range("L14").Formula = "=AVERAGE(G14:K14)"
$ws2.Range("L15:L38").Formula = "=AVERAGE(G15:K15)"

# ---------------------------------------------------------------------
# Sheet "500 rekordów"
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item(3)

$ws3.Range("L14").Formula = "=AVERAGE(G14:K14)"
$ws3.Range("L15:L38").Formula = "=AVERAGE(G15:K15)"

# ---------------------------------------------------------------------
# Restore each sheet's view (scroll position + selection), ending on
# sheet 3 so it remains the active tab.
# ---------------------------------------------------------------------
$ws1.Activate()
$excel.ActiveWindow.ScrollRow = 13
$excel.ActiveWindow.ScrollColumn = 4
$ws1.Range("P31").Select()

$ws2.Activate()
$excel.ActiveWindow.ScrollRow = 10
$excel.ActiveWindow.ScrollColumn = 1
$ws2.Range("O31").Select()

$ws3.Activate()
$excel.ActiveWindow.ScrollRow = 7
$excel.ActiveWindow.ScrollColumn = 1
$ws3.Range("R23").Select()
